$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowNum = 99
$values = @(98, 460, 32, 304, 26, 822, 23, 9, 32, 5203, 3421, 960, 40, 4243, 2, 38, 336, 736)

for ($i = 0; $i -lt $values.Length; $i++) {
    $col = $i + 1
    $ws.Cells.Item($rowNum, $col).Value = $values[$i]
}
